$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3090.974289607632
$ws.Range("D3").Value = 3090.97428960763
$ws.Range("D4").Value = 3090.97428960763

$ws.Range("D6").Value = 8831.355113164813
$ws.Range("D7").Value = 8831.355113164813

$ws.Range("D19").Value = 79354.61403309148
$ws.Range("D20").Value = 79354.61403309148

$ws.Range("D24").Value = 149626.5932405632
$ws.Range("D25").Value = 149626.593240563

$ws.Range("D28").Value = 2992.531864811288
$ws.Range("D29").Value = 2992.531864811288
$ws.Range("D30").Value = 149626.593240563

$ws.Range("D38").Value = -3390.728536028876
$ws.Range("D39").Value = -3390.728536028873

$ws.Range("D41").Value = 3390.728536028873

$ws.Range("D42").Value = 67814.57072057677
$ws.Range("D43").Value = 67814.57072057677
